$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = 44859
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("D13").Value = 44524
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 1000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1000
$ws.Range("D14").Value = 44460
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 2000
$ws.Range("P14").Value = 2000
$ws.Range("D15").Value = 44482
$ws.Range("J15").Value = 4000
$ws.Range("K15").Value = 900
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 950
$ws.Range("P15").Value = 950
$ws.Range("D16").Value = 44476
$ws.Range("J16").Value = 5000
$ws.Range("M16").Value = 1040
$ws.Range("P16").Value = 1040
$ws.Range("D17").Value = 44175
$ws.Range("J17").Value = 800
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1100
$ws.Range("M17").Value = 1050
$ws.Range("P17").Value = 1050
$ws.Range("D18").Value = 44515
$ws.Range("H18").Value = "Verde"
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 1200
$ws.Range("M18").Value = 1200
$ws.Range("P18").Value = 1200
$ws.Range("D19").Value = 44841
$ws.Range("J19").Value = 2500
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = 1000
$ws.Range("P19").Value = 1000
$ws.Range("D20").Value = 44830
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("P20").Value = 2000
$ws.Range("D21").Value = 44811
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 2500
$ws.Range("D22").Value = 44473
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = 1200
$ws.Range("O22").Value = "Provincia de Linares"
$ws.Range("P22").Value = 1200
$ws.Range("D23").Value = 44512
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 800
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = 800
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 800
$ws.Range("D24").Value = 44461
$ws.Range("J24").Value = 2500
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = 2000
$ws.Range("P24").Value = 2000
$ws.Range("D25").Value = 44468
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = 1500
$ws.Range("P25").Value = 1500
$ws.Range("D26").Value = 44477
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 1000
$ws.Range("P26").Value = 1000
$ws.Range("D27").Value = 44508
$ws.Range("J27").Value = 5000
$ws.Range("D28").Value = 44509
$ws.Range("J28").Value = 6000
$ws.Range("K28").Value = 800
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = 800
$ws.Range("P28").Value = 800
$ws.Range("D29").Value = 44474
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 1200
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = 1200
$ws.Range("N29").Value = "`$/kilo"
$ws.Range("P29").Value = 1200
$ws.Range("D30").Value = 44162
$ws.Range("H30").Value = "Verde"
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = 1000
$ws.Range("N30").Value = "`$/atado"
$ws.Range("P30").Value = 1000
$ws.Range("D31").Value = 44837
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1600
$ws.Range("L31").Value = 1600
$ws.Range("M31").Value = 1600
$ws.Range("P31").Value = 1600
$ws.Range("D32").Value = 44525
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = 1200
$ws.Range("P32").Value = 1200
$ws.Range("D33").Value = 44496
$ws.Range("H33").Value = "Verde"
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = 900
$ws.Range("P33").Value = 900
$ws.Range("D34").Value = 44827
$ws.Range("D35").Value = 44831
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 2000
$ws.Range("P35").Value = 2000
$ws.Range("D36").Value = 44160
$ws.Range("H36").Value = "Verde"
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 800
$ws.Range("L36").Value = 800
$ws.Range("M36").Value = 800
$ws.Range("P36").Value = 800
$ws.Range("D37").Value = 44826
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("J37").Value = 500
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = 2000
$ws.Range("P37").Value = 2000
$ws.Range("D38").Value = 44165
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 1200
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 1200
$ws.Range("O38").Value = "Provincia de Linares"
$ws.Range("P38").Value = 1200
$ws.Range("D39").Value = 44519
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 1100
$ws.Range("L39").Value = 1100
$ws.Range("M39").Value = 1100
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1100
$ws.Range("D40").Value = 44161
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = 1000
$ws.Range("P40").Value = 1000
$ws.Range("D41").Value = 44491
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 850
$ws.Range("L41").Value = 850
$ws.Range("M41").Value = 850
$ws.Range("O41").Value = "Provincia de Linares"
$ws.Range("P41").Value = 850
$ws.Range("D42").Value = 44516
$ws.Range("J42").Value = 3000
$ws.Range("K42").Value = 1000
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = 1000
$ws.Range("O42").Value = "Provincia de Limarí"
$ws.Range("P42").Value = 1000
$ws.Range("D43").Value = 44510
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 800
$ws.Range("M43").Value = 800
$ws.Range("P43").Value = 800
$ws.Range("D44").Value = 44469
$ws.Range("J44").Value = 3000
$ws.Range("K44").Value = 1200
$ws.Range("L44").Value = 1200
$ws.Range("M44").Value = 1200
$ws.Range("P44").Value = 1200
$ws.Range("D45").Value = 44505
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 800
$ws.Range("L45").Value = 800
$ws.Range("M45").Value = 800
$ws.Range("P45").Value = 800
$ws.Range("D46").Value = 44176
$ws.Range("H46").Value = "Verde"
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = 900
$ws.Range("P46").Value = 900
$ws.Range("D47").Value = 44848
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("J47").Value = 3000
$ws.Range("O47").Value = "Provincia de Linares"
$ws.Range("D48").Value = 44172
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 1000
$ws.Range("L48").Value = 1000
$ws.Range("M48").Value = 1000
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1000
$ws.Range("D49").Value = 44522
$ws.Range("H49").Value = "Verde"
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 1200
$ws.Range("L49").Value = 1200
$ws.Range("M49").Value = 1200
$ws.Range("P49").Value = 1200
$ws.Range("D50").Value = 44818
$ws.Range("J50").Value = 2000
$ws.Range("D51").Value = 44824
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("J51").Value = 1000
$ws.Range("K51").Value = 2800
$ws.Range("L51").Value = 2800
$ws.Range("M51").Value = 2800
$ws.Range("N51").Value = "`$/kilo"
$ws.Range("P51").Value = 2800
$ws.Range("Q51").Value = 1
$ws.Range("D52").Value = 44481
$ws.Range("J52").Value = 4000
$ws.Range("K52").Value = 900
$ws.Range("L52").Value = 900
$ws.Range("M52").Value = 900
$ws.Range("N52").Value = "`$/caja 10 kilos"
$ws.Range("P52").Value = 90
$ws.Range("Q52").Value = 10
$ws.Range("D53").Value = 44504
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 800
$ws.Range("L53").Value = 800
$ws.Range("M53").Value = 800
$ws.Range("P53").Value = 800
$ws.Range("D54").Value = 44166
$ws.Range("H54").Value = "Verde"
$ws.Range("J54").Value = 1500
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 1000
$ws.Range("M54").Value = 1000
$ws.Range("P54").Value = 1000
$ws.Range("D55").Value = 44855
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 1200
$ws.Range("L55").Value = 1200
$ws.Range("M55").Value = 1200
$ws.Range("P55").Value = 1200
$ws.Range("D56").Value = 44494
$ws.Range("J56").Value = 4000
$ws.Range("K56").Value = 900
$ws.Range("L56").Value = 900
$ws.Range("M56").Value = 900
$ws.Range("P56").Value = 900
$ws.Range("D57").Value = 44518
$ws.Range("J57").Value = 3000
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 1000
$ws.Range("M57").Value = 1000
$ws.Range("P57").Value = 1000
$ws.Range("D58").Value = 44488
$ws.Range("H58").Value = "Verde"
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 850
$ws.Range("L58").Value = 900
$ws.Range("M58").Value = 875
$ws.Range("P58").Value = 875
$ws.Range("D59").Value = 44816
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("J59").Value = 500
$ws.Range("K59").Value = 2800
$ws.Range("L59").Value = 2800
$ws.Range("M59").Value = 2800
$ws.Range("P59").Value = 2800
$ws.Range("D60").Value = 44168
$ws.Range("D61").Value = 44169
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = 1000
$ws.Range("P61").Value = 1000
$ws.Range("D62").Value = 44517
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 1100
$ws.Range("L62").Value = 1100
$ws.Range("M62").Value = 1100
$ws.Range("P62").Value = 1100
$ws.Range("D63").Value = 44463
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = 2000
$ws.Range("P63").Value = 2000
$ws.Range("D64").Value = 44487
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 800
$ws.Range("L64").Value = 800
$ws.Range("M64").Value = 800
$ws.Range("O64").Value = "Provincia de Linares"
$ws.Range("P64").Value = 800
$ws.Range("D65").Value = 44167
$ws.Range("K65").Value = 1000
$ws.Range("L65").Value = 1000
$ws.Range("M65").Value = 1000
$ws.Range("O65").Value = "Región del Maule"
$ws.Range("P65").Value = 1000
$ws.Range("D66").Value = 44462
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 1800
$ws.Range("L66").Value = 2000
$ws.Range("M66").Value = 1900
$ws.Range("P66").Value = 1900
$ws.Range("D67").Value = 44498
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 900
$ws.Range("L67").Value = 900
$ws.Range("M67").Value = 900
$ws.Range("O67").Value = "Provincia de Linares"
$ws.Range("P67").Value = 900
$ws.Range("D68").Value = 44455
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = 2400
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 2400
$ws.Range("D69").Value = 44858
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("K69").Value = 1000
$ws.Range("L69").Value = 1100
$ws.Range("M69").Value = 1050
$ws.Range("O69").Value = "Provincia de Linares"
$ws.Range("P69").Value = 1050
$ws.Range("D70").Value = 44484
$ws.Range("L70").Value = 900
$ws.Range("M70").Value = 840
$ws.Range("P70").Value = 840
$ws.Range("D71").Value = 44495
$ws.Range("H71").Value = "Verde"
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 900
$ws.Range("L71").Value = 900
$ws.Range("M71").Value = 900
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 900
$ws.Range("D72").Value = 44503
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 800
$ws.Range("L72").Value = 800
$ws.Range("M72").Value = 800
$ws.Range("O72").Value = "Provincia de Linares"
$ws.Range("P72").Value = 800
$ws.Range("D73").Value = 44832
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = 2000
$ws.Range("P73").Value = 2000
$ws.Range("D74").Value = 44490
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 850
$ws.Range("M74").Value = 875
$ws.Range("O74").Value = "Región del Maule"
$ws.Range("P74").Value = 875
$ws.Range("D75").Value = 44845
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 1200
$ws.Range("L75").Value = 1200
$ws.Range("M75").Value = 1200
$ws.Range("P75").Value = 1200
$ws.Range("A76").Value = 5
$ws.Range("B76").Value = "Macroferia Regional de Talca"
$ws.Range("C76").Value = "Maule"
$ws.Range("D76").Value = 44497
$ws.Range("D76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E76").Value = 7
$ws.Range("F76").Value = 300000000
$ws.Range("G76").Value = "Espárragos"
$ws.Range("H76").Value = "Verde"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 900
$ws.Range("L76").Value = 900
$ws.Range("M76").Value = 900
$ws.Range("N76").Value = "`$/kilo"
$ws.Range("O76").Value = "Provincia de Linares"
$ws.Range("P76").Value = 900
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"
$ws.Range("A77").Value = 5
$ws.Range("B77").Value = "Macroferia Regional de Talca"
$ws.Range("C77").Value = "Maule"
$ws.Range("D77").Value = 44489
$ws.Range("D77").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E77").Value = 7
$ws.Range("F77").Value = 300000000
$ws.Range("G77").Value = "Espárragos"
$ws.Range("H77").Value = "Verde"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 900
$ws.Range("L77").Value = 900
$ws.Range("M77").Value = 900
$ws.Range("N77").Value = "`$/kilo"
$ws.Range("O77").Value = "Provincia de Linares"
$ws.Range("P77").Value = 900
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"
